# Upload new version with timestamp
# Populate the pharmacy "transactions count" report with the product rows,
# move the footer/total placeholders down to their new positions, and
# refresh the report timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Preserve the two "placeholder" styled regions before we touch them:
#    - the empty totals strip currently at K5:N5 -> becomes the new
#      totals row at K18:N18
#    - the footer strip currently at row 6 -> becomes the new footer
#      row at row 19
# ---------------------------------------------------------------------
$ws.Range("K5:N5").Copy()
$ws.Range("K18:N18").PasteSpecial(-4122)

$ws.Range("A6:N6").Copy()
$ws.Range("A19:N19").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Remove the old merges for those two placeholder regions so the
#    rows can be reused for real product data.
# ---------------------------------------------------------------------
$ws.Range("K5:N5").UnMerge()
$ws.Range("A6:E6").UnMerge()
$ws.Range("F6:G6").UnMerge()
$ws.Range("I6:N6").UnMerge()

# ---------------------------------------------------------------------
# 3) Build the first product row (row 4). The "in:out" counter column
#    (H) and the product-name column (B) must be stored as literal
#    text (not re-interpreted as a time/ratio/number), so switch those
#    ranges to Text format before writing the values.
# ---------------------------------------------------------------------
$ws.Range("B4:G4").NumberFormat = "@"
$ws.Range("H4:K4").NumberFormat = "@"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ANGIOFOX (EFFOX) 25MG LONG 30 CAPS."
$ws.Range("H4").Value = "0:0"
$ws.Range("L4").Value = 114
$ws.Range("N4").Value = 1

# ---------------------------------------------------------------------
# 4) Stamp that same formatting down through row 17, then merge each
#    row's column groups the same way row 4 already is.
# ---------------------------------------------------------------------
$ws.Range("A4:N4").Copy()
$ws.Range("A5:N17").PasteSpecial(-4122)

$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 24.75
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 24.75
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 25.5
$ws.Rows.Item(19).RowHeight = 17.25

$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()
$ws.Range("B6:G6").Merge()
$ws.Range("H6:K6").Merge()
$ws.Range("L6:M6").Merge()
$ws.Range("B7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("B9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("B10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("B11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("B12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("B13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("B14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("B16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("B17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()

$ws.Range("K18:N18").Merge()
$ws.Range("A19:E19").Merge()
$ws.Range("F19:G19").Merge()
$ws.Range("I19:N19").Merge()

# ---------------------------------------------------------------------
# 5) Fill in the remaining product rows (5 through 17).
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "AUGMENTIN 457MG/5ML SUSP. 70 ML"
$ws.Range("H5").Value = "1:0"
$ws.Range("L5").Value = 137
$ws.Range("N5").Value = 1

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "BLOKATENS 10/160MG 28 F.C.TABS."
$ws.Range("H6").Value = "0:0"
$ws.Range("L6").Value = 160
$ws.Range("N6").Value = 1

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "COLOVATIL 30 F.C. TABS"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = 63
$ws.Range("N7").Value = 1

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "GAVISCON LIQUID 24 SACHETS 10 ML"
$ws.Range("H8").Value = "0:9"
$ws.Range("L8").Value = 12
$ws.Range("N8").Value = 0.04

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "GINKGO BILOBA 30 CAPS."
$ws.Range("H9").Value = "0:0"
$ws.Range("L9").Value = 186
$ws.Range("N9").Value = 1

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "MILGA ADVANCE 30 F.C. TABS"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 136.5
$ws.Range("N10").Value = 1

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "PERLOC 40MG 14 F.C.TAB."
$ws.Range("H11").Value = "0:0"
$ws.Range("L11").Value = 68.25
$ws.Range("N11").Value = 1

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "RHINEX 0.05% INFANTILE NASAL DROPS 10 ML"
$ws.Range("H12").Value = "2:0"
$ws.Range("L12").Value = 18
$ws.Range("N12").Value = 1

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "RIVO 320MG 20*10 TABS"
$ws.Range("H13").Value = "1:2"
$ws.Range("L13").Value = 14.1
$ws.Range("N13").Value = 0.1

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "VASTAREL MR 35MG 30 F.C.TAB."
$ws.Range("H14").Value = "2:0"
$ws.Range("L14").Value = 175
$ws.Range("N14").Value = 1

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "WATER FOR INJECTION AMP. 5 ML"
$ws.Range("H15").Value = "7816:0"
$ws.Range("L15").Value = 2.5
$ws.Range("N15").Value = 1

$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "سويت كوكو"
$ws.Range("H16").Value = "22:0"
$ws.Range("L16").Value = 25
$ws.Range("N16").Value = 1

$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "مرطب شفاه لونا جوز هند ابيض"
$ws.Range("H17").Value = "3:0"
$ws.Range("L17").Value = 20
$ws.Range("N17").Value = 1

# ---------------------------------------------------------------------
# 6) Grand total row.
# ---------------------------------------------------------------------
$ws.Range("K18").Value = 1131.35

# ---------------------------------------------------------------------
# 7) Refresh the footer timestamp / page / credit line (same values,
#    just relocated from row 6 to row 19).
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Monday, 5 January, 2026 9:53 AM"
$ws.Range("F19").Value = "1/1"
$ws.Range("I19").Value = "developed by : Abdelaziz Talaat"

Write-Output "Report rows populated"
